$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Sheet1")

# --- Update score values on Sheet1 ---
$ws.Range("C4").Value = 85

$ws.Range("C5").Value = 72
$ws.Range("D5").Value = 66

$ws.Range("C7").Value = 74

$ws.Range("C17").Value = 35

$ws.Range("C19").Value = 48

$ws.Range("C22").Value = 35
$ws.Range("D22").Value = 26

$ws.Range("C25").Value = 36
$ws.Range("D25").Value = 31

$ws.Range("C27").Value = 59
$ws.Range("D27").Value = 44

$ws.Range("D28").Value = 44

$ws.Range("C33").Value = 51

$ws.Range("D38").Value = 54

$ws.Range("C40").Value = 80
$ws.Range("D40").Value = 70

$ws.Range("C41").Value = 61

$ws.Range("D44").Value = 57

$ws.Range("D46").Value = 52

$ws.Range("D53").Value = 75

$ws.Range("C54").Value = 54
$ws.Range("D54").Value = 33

$ws.Range("C61").Value = 38
$ws.Range("D61").Value = 32

$ws.Range("C68").Value = 58
$ws.Range("D68").Value = 45

$ws.Range("C73").Value = 58
$ws.Range("D73").Value = 40

$ws.Range("C77").Value = 133
$ws.Range("D77").Value = 132

$ws.Range("C84").Value = 190
$ws.Range("D84").Value = 125

$ws.Range("C92").Value = 245
$ws.Range("D92").Value = 175

# Grand total row
$ws.Range("C93").Value = 5450
$ws.Range("D93").Value = 4377

# Select the full first row, as stored in the updated sheet view
$ws.Range("A1:XFD1").Select()

# --- Remove the now-unused extra worksheets ---
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()
